{"js": "// 1) \"rif.: allegato 6 del Protocollo eGLU:\" used to be split across several\n//    runs (artifacts of Word's spell-checker, each wrapped in <w:proofErr/>\n//    start/end pairs). Re-typing the same text over that range collapses it\n//    back into a single run (same visible text, same run formatting) the\n//    way the target document has it.\nconst intro = context.document.body.search(\n  \"rif.: allegato 6 del Protocollo eGLU:\",\n  { matchCase: true }\n);\nintro.load(\"items\");\nawait context.sync();\n\nif (intro.items.length > 0) {\n  intro.items[0].insertText(\n    \"rif.: allegato 6 del Protocollo eGLU:\",\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 2) The \"Valutazione totale\" score in the results table changed from 75 to 80.\nconst score = context.document.body.search(\"75\", {\n  matchCase: true,\n  matchWholeWord: true,\n});\nscore.load(\"items\");\nawait context.sync();\n\nscore.items.forEach((item) => {\n  item.insertText(\"80\", Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"rif.: allegato 6 del Protocollo eGLU:\" used to be split across several\n#    runs (artifacts of Word's spell-checker, each wrapped in a\n#    proofErr start/end pair). Re-typing the same text over that range\n#    collapses it back into a single run (same visible text, same run\n#    formatting), matching the target document.\n$introRange = $d.Content\n$introRange.Find.Execute(\n    \"rif.: allegato 6 del Protocollo eGLU:\",\n    $true, $false, $false, $false, $false,\n    $true, 1, $false,\n    \"rif.: allegato 6 del Protocollo eGLU:\",\n    2\n) | Out-Null\n\n# 2) The \"Valutazione totale\" score in the results table changed from 75 to 80.\n$scoreRange = $d.Content\n$scoreRange.Find.Execute(\n    \"75\",\n    $true, $true, $false, $false, $false,\n    $true, 1, $false,\n    \"80\",\n    2\n) | Out-Null\n"}
